# Clear the unused "2nd run" header labels (Time_2, Car_2, Track_2 ... Track_5)
# while keeping their formatting/style.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G1:R1").ClearContents()

# Clear the Car_2 (H) and Track_2 (I) data columns for all player rows (2-38):
# these were defaulted to "Koi" / "Mt. Hairpin" and should become blank (0/NA -> empty).
$ws.Range("H2:I38").ClearContents()

# Update the active selection to match the edited file (J12).
$ws.Range("J12").Select()
